# Commit: Création des scripts de données pour les styles et quelques paramètres
#
# 1) Rename "Sheet1" -> "Tables"
# 2) Populate the "Données" sheet with the SQL seed-script reference tables
#    (Style, TypeParam) by re-using formatting already present on "Tables"
#    via copy / paste-special (formats only), and a new bold+underlined
#    font for the sub-section headers.
# 3) Restore the selections / column width / page setup that Excel leaves
#    behind after this kind of editing session.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "Tables"

$ws2 = $wb.Worksheets.Item("Données")

$ws2.Range("B1").Value = "Style"
$ws1.Range("B2").Copy()
$ws2.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("B2").Value = "Id"
$ws1.Range("B3").Copy()
$ws2.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("B3").Value = "Nom"
$ws1.Range("D28").Copy()
$ws2.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("B5").Value = "insert into style (Nom) values ('Rock')"
$ws2.Range("B6").Value = "insert into style (Nom) values ('Classique')"
$ws2.Range("B7").Value = "insert into style (Nom) values ('Urbain')"
$ws2.Range("B8").Value = "insert into style (Nom) values ('Sport')"
$ws2.Range("B9").Value = "insert into style (Nom) values ('Retro')"
$ws2.Range("B10").Value = "insert into style (Nom) values ('Traditionnel')"
$ws2.Range("B13").Value = "TypeParam"
$ws1.Range("B2").Copy()
$ws2.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("B14").Value = "Id"
$ws1.Range("B3").Copy()
$ws2.Range("B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("B15").Value = "TypeId"
$ws1.Range("D19").Copy()
$ws2.Range("B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("B16").Value = "TypeLib"
$ws1.Range("D19").Copy()
$ws2.Range("B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("B17").Value = "ParamCode"
$ws1.Range("D19").Copy()
$ws2.Range("B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("B18").Value = "ParamLib"
$ws1.Range("D28").Copy()
$ws2.Range("B18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("B20").Value = "Type statut personne (historique)"
$ws2.Range("B20").Font.Bold = $true
$ws2.Range("B20").Font.Underline = $true
$ws2.Range("B21").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (1,'PersonneStatus','Abonne','Personne abonnée')"
$ws2.Range("B25").Value = "Type de photo"
$ws2.Range("B25").Font.Bold = $true
$ws2.Range("B25").Font.Underline = $true
$ws2.Range("B22").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (1,'PersonneStatus','EnAttente','Personne en attente')"
$ws2.Range("B23").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (1,'PersonneStatus','Conseiller','Personne conseiller')"
$ws2.Range("B26").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (2,'PhotoType','Habille','Photo de personne habillé')"
$ws2.Range("B27").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (2,'PhotoType','Vetement','Photo de vêtement')"
$ws2.Range("B29").Value = "Type de message"
$ws2.Range("B29").Font.Bold = $true
$ws2.Range("B29").Font.Underline = $true
$ws2.Range("B30").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (3,'MessageType','Message','Message simple')"
$ws2.Range("B31").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (3,'MessageType','Notification','Notification')"
$ws2.Range("B32").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (3,'MessageType','Information','Information')"
$ws2.Range("B34").Value = "Type statut demande"
$ws2.Range("B34").Font.Bold = $true
$ws2.Range("B34").Font.Underline = $true
$ws2.Range("B35").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','AttenteDemandeur','AttenteDemandeur')"
$ws2.Range("B36").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','AttenteConseiller','AttenteConseiller')"
$ws2.Range("B37").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','Accepte','Accepte')"
$ws2.Range("B38").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','RefusDemandeur','RefusDemandeur')"
$ws2.Range("B39").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','RefusConseiller','RefusConseiller')"
$ws2.Range("B40").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','AnnulAdmin','AnnulAdmin')"
$ws2.Range("B41").Value = "insert into typeparam (TypeId, TypeLib, ParamCode, ParamLib) values (4,'DemandeStatus','Termine','Termine')"
$ws2.Range("B43").Value = "Type de vêtement"
$ws2.Range("B43").Font.Bold = $true
$ws2.Range("B43").Font.Underline = $true
$ws2.Range("E43").Value = "VetementType"
$ws2.Range("B44").Value = "Tete"
$ws2.Range("B50").Value = "Type compteur"
$ws2.Range("B50").Font.Bold = $true
$ws2.Range("B50").Font.Underline = $true
$ws2.Range("E50").Value = "CompteurType"
$ws2.Range("B51").Value = "Abonne"
$ws2.Range("B52").Value = "Conseiller"
$ws2.Range("B54").Value = "Type de notification"
$ws2.Range("B54").Font.Bold = $true
$ws2.Range("B54").Font.Underline = $true
$ws2.Range("E54").Value = "NotifType"
$ws2.Range("B56").Value = "PropositionAccept,"
$ws2.Range("B57").Value = "PropositionReject,"
$ws2.Range("B58").Value = "PropositionCreation"
$ws2.Range("B59").Value = "DemandAccept,"
$ws2.Range("B60").Value = "DemandReject"
$ws2.Range("B45").Value = "Buste,"
$ws2.Range("B46").Value = "Jambe,"
$ws2.Range("B47").Value = "Pied,"
$ws2.Range("B48").Value = "Accessoire"
$ws2.Range("B55").Value = "DemandCreation"

# Column B on "Données" is a bit wider than the default so the long
# "insert into ..." strings are easier to read.
$ws2.Columns.Item(2).ColumnWidth = 10.2

# Match the print setup that was configured for this sheet.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Restore the selections left on each sheet.
$ws1.Range("J15:J20").Select()
$ws2.Range("D47").Select()
